$d = $word.ActiveDocument

# Locate the paragraph that reads exactly "BASE DE DATOS" and is immediately
# followed by the "Adjuntamos un archivo Excel con todas las tablas" paragraph
# (there is an earlier, unrelated "BASE DE DATOS" bullet higher up in the doc,
# so we must disambiguate using the following paragraph's text).
$target = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "BASE DE DATOS`r" -and $i -lt $d.Paragraphs.Count) {
        $nextPara = $d.Paragraphs.Item($i + 1)
        if ($nextPara.Range.Text -eq "Adjuntamos un archivo Excel con todas las tablas`r") {
            $target = $i
            break
        }
    }
}

if ($target -eq 0) {
    throw "Could not find the target 'BASE DE DATOS' paragraph"
}

# Sanity-check the two blank paragraphs that must follow "Adjuntamos..." before
# we remove anything, so we never delete unexpected content.
$blank1 = $d.Paragraphs.Item($target + 2).Range.Text
$blank2 = $d.Paragraphs.Item($target + 3).Range.Text
if ($blank1 -ne "`r" -or $blank2 -ne "`r") {
    throw "Unexpected content after the 'Adjuntamos...' paragraph; aborting"
}

# Build the new paragraph text ("BASE DE DATOS" + ": " + "ver carpeta BASE DE
# DATOS") as three separate runs, inserted as leading content of the target
# paragraph so the paragraph keeps its own properties (list style/numbering).
$targetPara = $d.Paragraphs.Item($target)
$insertionPoint = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)
$newRunsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:r><w:t>BASE DE DATOS</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">: </w:t></w:r>' + `
    '<w:r><w:t>ver carpeta BASE DE DATOS</w:t></w:r>' + `
    '</w:p>'
[void]$insertionPoint.InsertXML($newRunsXml)

# Remove the now-duplicated original "BASE DE DATOS" run that got pushed to
# the end of the paragraph.
$mergedPara = $d.Paragraphs.Item($target)
$markerPos = $mergedPara.Range.End - 1
$oldRunLength = "BASE DE DATOS".Length
$dupRange = $d.Range($markerPos - $oldRunLength, $markerPos)
$dupRange.Delete()

# Delete the following three paragraphs entirely: the "Adjuntamos un archivo
# Excel con todas las tablas" paragraph and the two blank paragraphs after it.
$deleteStart = $d.Paragraphs.Item($target + 1).Range.Start
$deleteEnd = $d.Paragraphs.Item($target + 3).Range.End
$d.Range($deleteStart, $deleteEnd).Delete()
